{"js": "// Add a new row \"55.7 Ordenaci\u00f3n (ascendente y descendente) de modelos por campos\"\n// right after the existing \"55.6 Manual de usuario y de administrador\" row, in the\n// single table that lists the \"Codigos de Recopilacion de Esfuerzos\".\n\n// Locate the anchor row by searching for the text of the row that precedes the\n// new one - this avoids hard-coding a row index that could drift if the table\n// changes shape upstream.\nconst searchResults = context.document.body.search(\n  \"Manual de usuario y de administrador\",\n  { matchCase: false, matchWholeWord: false }\n);\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find the anchor text \"Manual de usuario y de administrador\".');\n}\n\nconst anchorCell = searchResults.items[0].parentTableCell;\nconst anchorRow = anchorCell.parentRow;\n\n// Insert a brand-new row right after the anchor row and seed its two cells in\n// one shot; the host mirrors the surrounding row's cell widths/alignment and\n// the paragraph/run formatting (Calibri, black, 18 half-points, es-ES) onto\n// the freshly inserted cells automatically.\nanchorRow.insertRows(\"After\", 1, [\n  [\"55.7\", \"Ordenaci\u00f3n (ascendente y descendente) de modelos por campos\"]\n]);\n\nawait context.sync();\n", "ps1": "# Add a new row \"55.7 Ordenaci\u00f3n (ascendente y descendente) de modelos por campos\"\n# right after the existing \"55.6 Manual de usuario y de administrador\" row, in the\n# single table that lists the \"Codigos de Recopilacion de Esfuerzos\".\n\n$d = $word.ActiveDocument\n\n# Locate the anchor row by searching for the text of the row that precedes the\n# new one - this avoids hard-coding a row index that could drift if the table\n# changes shape upstream.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Manual de usuario y de administrador\")\nif (-not $found) {\n    throw \"Could not find the anchor text 'Manual de usuario y de administrador'.\"\n}\n\n$anchorCell = $rng.Cells.Item(1)\n$table = $anchorCell.Tables.Item(1)\n$anchorRowIndex = $anchorCell.RowIndex\n\n# Row that currently sits right after the anchor row; passing it as the\n# \"before\" row to Rows.Add makes the brand-new row land exactly between the\n# two, i.e. right after the anchor row.\n$followingRow = $table.Rows.Item($anchorRowIndex + 1)\n$newRow = $table.Rows.Add($followingRow)\n\n$newRow.Cells.Item(1).Range.Text = \"55.7\"\n$newRow.Cells.Item(2).Range.Text = \"Ordenaci\u00f3n (ascendente y descendente) de modelos por campos\"\n"}
